# feat: add 2022-Q3 data
#
# Inserts a new worksheet "2022-Q3" (fund-holding detail, same shape as the
# existing "2021-Q1" / "2020-Q4" sheets) positioned right after "总计" and
# before "2021-Q1", and adds the corresponding summary row to "总计".
#
# NOTE: worksheet references are re-fetched by name AFTER every call that
# inserts/moves a sheet (Worksheets.Add shifts tab positions, which can
# stale out previously-grabbed COM handles for sheets after the insertion
# point) -- so every $wb.Worksheets.Item(...) lookup below happens only
# after the sheet collection has reached its final shape.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right before "2021-Q1" so the tab
#    order becomes: 总计, 2022-Q3, 2021-Q1, 2020-Q4
# ---------------------------------------------------------------------
$new = $wb.Worksheets.Add($wb.Worksheets.Item("2021-Q1"))
$new.Name = "2022-Q3"

# Re-fetch every sheet handle now that the collection is in its final shape.
$zj = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")
$q1 = $wb.Worksheets.Item("2021-Q1")
$q4 = $wb.Worksheets.Item("2020-Q4")

$q3.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$q3.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$q3.PageSetup.TopMargin = $excel.InchesToPoints(1)
$q3.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$q3.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$q3.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)

# ---------------------------------------------------------------------
# 2. Header row — same look as the other fund-holding sheets (bold,
#    centered, bordered => style copied from "2021-Q1" header).
# ---------------------------------------------------------------------
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q1.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Data rows (fund code / numeric-looking figures are stored as TEXT
#    in this workbook, so force "@" text format before assigning them,
#    otherwise Excel would silently coerce "001556" -> 1556, "25.09" ->
#    a real number, etc.). The rank/index columns (A, H) are real
#    numbers.
# ---------------------------------------------------------------------
function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$rows = @(
    @{ a = 0; b = "001556"; c = "天弘中证500指数增强A";             d = "25.09"; e = "94.15"; f = "1.44"; g = "0.3613"; h = 9 },
    @{ a = 1; b = "001557"; c = "天弘中证500指数增强C";             d = "12.94"; e = "94.15"; f = "1.44"; g = "0.1863"; h = 9 },
    @{ a = 2; b = "000270"; c = "建信灵活配置混合";                   d = "2.27";  e = "94.21"; f = "1.01"; g = "0.0229"; h = 2 },
    @{ a = 3; b = "003242"; c = "创金合信量化发现灵活配置混合C"; d = "0.40";  e = "92.08"; f = "1.49"; g = "0.0060"; h = 9 },
    @{ a = 4; b = "003241"; c = "创金合信量化发现灵活配置混合A"; d = "0.32";  e = "92.08"; f = "1.49"; g = "0.0048"; h = 9 }
)

$r = 2
foreach ($row in $rows) {
    $q3.Range("A$r").Value = $row.a
    Set-TextCell $q3.Range("B$r") $row.b
    $q3.Range("C$r").Value = $row.c
    Set-TextCell $q3.Range("D$r") $row.d
    Set-TextCell $q3.Range("E$r") $row.e
    Set-TextCell $q3.Range("F$r") $row.f
    Set-TextCell $q3.Range("G$r") $row.g
    $q3.Range("H$r").Value = $row.h
    $r = $r + 1
}

# Column A on the data rows carries the same bordered/centered style as
# the header (matches the other fund sheets).
$q1.Range("A2").Copy()
$q3.Range("A2:A6").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4. "总计" (summary) sheet — shift the existing two rows down one slot
#    and insert the new 2022-Q3 summary row above them.
# ---------------------------------------------------------------------
$oldRow2B = $zj.Range("B2").Value()
$oldRow2C = $zj.Range("C2").Value()
$oldRow2D = $zj.Range("D2").Value()
$oldRow3B = $zj.Range("B3").Value()
$oldRow3C = $zj.Range("C3").Value()
$oldRow3D = $zj.Range("D3").Value()

$zj.Range("A4").Value = 2
$zj.Range("B4").Value = $oldRow3B
$zj.Range("C4").Value = $oldRow3C
$zj.Range("D4").Value = $oldRow3D

$zj.Range("A3").Value = 1
$zj.Range("B3").Value = $oldRow2B
$zj.Range("C3").Value = $oldRow2C
$zj.Range("D3").Value = $oldRow2D

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q3"
$zj.Range("C2").Value = 5
$zj.Range("D2").Value = 0.58

$zj.Range("A2").Copy()
$zj.Range("A3:A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 5. Restore the originally-active tab ("2020-Q4").
# ---------------------------------------------------------------------
$q4.Activate()

Write-Output "2022-Q3 sheet added"
